$wb = $excel.ActiveWorkbook

# --- Sheet "Mensal" (monthly data): append row 14 ---
$wsMensal = $wb.Worksheets.Item("Mensal")

# Copy formatting (date style) from the last existing row so the new
# date cell keeps the same style index as the rest of column A.
$wsMensal.Range("A13").Copy()
$wsMensal.Range("A14").PasteSpecial(-4122)

$wsMensal.Range("A14").Value = 44150
$wsMensal.Range("B14").Value = 154.96
$wsMensal.Range("C14").Value = 175.21
$wsMensal.Range("D14").Value = -11.55

# --- Sheet "Diario" (daily data): append rows 368-382 ---
$wsDiario = $wb.Worksheets.Item("Diario")

$diarioData = @(
    @(44136, 156.12, 175.21, -10.89),
    @(44137, 155.39, 175.21, -11.31),
    @(44138, 155.17, 175.21, -11.43),
    @(44139, 154.85, 175.21, -11.62),
    @(44140, 154.5,  175.21, -11.82),
    @(44141, 154.33, 175.21, -11.91),
    @(44142, 154.03, 175.21, -12.09),
    @(44143, 154.28, 175.21, -11.94),
    @(44144, 153.95, 175.21, -12.13),
    @(44145, 154.26, 175.21, -11.96),
    @(44146, 155.22, 175.21, -11.41),
    @(44147, 155.18, 175.21, -11.43),
    @(44148, 156.41, 175.21, -10.73),
    @(44149, 154.66, 175.21, -11.73),
    @(44150, 156.12, 175.21, -10.89)
)

$startRow = 368
$lastExistingRow = $startRow - 1

for ($i = 0; $i -lt $diarioData.Count; $i++) {
    $r = $startRow + $i
    $row = $diarioData[$i]

    # Copy formatting (date style) from the last existing row's date
    # cell so every new date cell keeps the same style index.
    $wsDiario.Range("A$lastExistingRow").Copy()
    $wsDiario.Range("A$r").PasteSpecial(-4122)

    $wsDiario.Cells.Item($r, 1).Value = $row[0]
    $wsDiario.Cells.Item($r, 2).Value = $row[1]
    $wsDiario.Cells.Item($r, 3).Value = $row[2]
    $wsDiario.Cells.Item($r, 4).Value = $row[3]
}
